$wb = $excel.ActiveWorkbook

# --- Rename existing sheets / add the new "API responses" sheet ---
$messagesSheet = $wb.Worksheets.Item("Blad1")
$messagesSheet.Name = "API messages"

# Insert the new sheet right after "API messages" (so final order is
# API messages, API responses, Blad2)
$responsesSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $messagesSheet)
$responsesSheet.Name = "API responses"

# --- Fill in the new sheet's data ---
$responsesSheet.Range("A1").Value = "code"
$responsesSheet.Range("B1").Value = "betekenis"

$responsesSheet.Range("A2").Value = 0
$responsesSheet.Range("B2").Value = "item succesvol gecreëerd"

$responsesSheet.Range("A3").Value = 1
$responsesSheet.Range("B3").Value = "operatie failed, databasefout"

$responsesSheet.Range("A4").Value = 2
$responsesSheet.Range("B4").Value = "operatie failed, item bestaat al in DB met zelfde of hoger versienummer"

# --- Formatting: bold header row with a bottom border ---
$headerRange = $responsesSheet.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Borders.Item(9).LineStyle = 1
$headerRange.Borders.Item(9).Weight = 2

# "betekenis" header left-aligned
$responsesSheet.Range("B1").HorizontalAlignment = -4131

# Numeric code column centered
$responsesSheet.Range("A2:A4").HorizontalAlignment = -4108

# --- Column widths ---
$responsesSheet.Columns.Item(2).AutoFit() | Out-Null

# --- Page setup to mirror the other worksheet ---
$responsesSheet.PageSetup.PaperSize = 9
$responsesSheet.PageSetup.Orientation = 1

# --- Selection / active cell state for the new sheet ---
$responsesSheet.Range("B7").Select() | Out-Null

# --- Workbook window / file metadata ---
$wb.Windows.Item(1).Left = 768
$wb.Windows.Item(1).Top = 768
$wb.Windows.Item(1).Width = 17280
$wb.Windows.Item(1).Height = 8964
